$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the three new worksheets after Sheet1, in order: Sheet2, shitttt, shitttt2
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "shitttt"

$ws4 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws3)
$ws4.Name = "shitttt2"

# Fill remaining rows on Sheet1 (feedback/bug entries + blank separators)
$ws1.Range("A2").Value = 'feedback'
$ws1.Range("B2").Value = 'Amazing games'
$ws1.Range("C2").Value = 'positive'
$ws1.Range("D2").Value = 'What a beauty!!!! Amazing games - both graphics, content, and feedback that children receive! Yesterday I tried with my 5-year-old daughter, it was difficult for her, she had to make an effort, but the satisfaction after the climb was worth the effort :).'
$ws1.Range("A3").Value = 'bug'
$ws1.Range("B3").Value = 'Navigation game issue'
$ws1.Range("C3").Value = 'negative'
$ws1.Range("D3").Value = 'In the navigation game, it was a bit disturbing for her and also for the patients I tried with today to collect minerals when there is a note about maintaining an upright back - the correct posture window hides the path and is a bit stressful. Maybe after the initial guidance, only voice guidance would be sufficient?'
$ws1.Range("A4").Value = 'feedback'
$ws1.Range("B4").Value = 'Amazing games with graphics, content, and feedback'
$ws1.Range("C4").Value = 'positive'
$ws1.Range("D4").Value = 'What a beauty!!!! Amazing games - both graphics, content, and feedback that children receive! Yesterday I tried with my 5-year-old daughter, it was hard for her, she had to make an effort, but the satisfaction after the climb was worth the effort :).'
$ws1.Range("A5").Value = 'bug'
$ws1.Range("B5").Value = 'Issue with navigation game and collecting minerals'
$ws1.Range("C5").Value = 'negative'
$ws1.Range("D5").Value = 'In the navigation game, it was a bit disturbing for her and also for the patients I tried with today to collect the minerals when a note appears to maintain a straight back - the correct posture window hides the path and is a bit stressful. Maybe after the first guidance, only a voice guidance could be sufficient?'
$ws1.Range("A6").Value = 'feedback'
$ws1.Range("B6").Value = 'Amazing games and graphics'
$ws1.Range("C6").Value = 'positive'
$ws1.Range("D6").Value = 'What a beauty! Amazing games - graphics, content, and feedback for children! I tried it yesterday with my 5-year-old daughter, it was hard for her, she had to make an effort, but the satisfaction after the climb was worth the effort.'
$ws1.Range("A7").Value = 'bug'
$ws1.Range("B7").Value = 'Navigation game issue'
$ws1.Range("C7").Value = 'negative'
$ws1.Range("D7").Value = 'In the navigation game, it was a bit disturbing for her and also for the patients I tried with today to collect the minerals when there is a note about keeping a straight back - the proper posture window hides the route and is a bit stressful. Maybe after the first guidance, a voice direction would be enough?'
$ws1.Range("A8").Value = ' '
$ws1.Range("A9").Value = 'feedback'
$ws1.Range("B9").Value = 'Amazing games and graphics'
$ws1.Range("C9").Value = 'positive'
$ws1.Range("D9").Value = 'What a beauty! Amazing games - graphics, content, and feedback for children! I tried it yesterday with my 5-year-old daughter, it was hard for her, she had to make an effort, but the satisfaction after the climb was worth the effort.'
$ws1.Range("A10").Value = 'bug'
$ws1.Range("B10").Value = 'Navigation game issue'
$ws1.Range("C10").Value = 'negative'
$ws1.Range("D10").Value = 'In the navigation game, it was a bit disturbing for her and also for the patients I tried with today to collect the minerals when there is a note about keeping a straight back - the proper posture window hides the route and is a bit stressful. Maybe after the first guidance, a voice direction would be enough?'
$ws1.Range("A11").Value = 'feedback'
$ws1.Range("B11").Value = 'Amazing games and positive experience'
$ws1.Range("C11").Value = 'positive'
$ws1.Range("D11").Value = 'What a beauty!!!! Amazing games - both graphics, content and feedback that children receive! Yesterday I tried with my 5-year-old daughter, it was hard for her, she had to make an effort, but the satisfaction after the climb was worth the effort :).'
$ws1.Range("A12").Value = 'bug'
$ws1.Range("B12").Value = 'Issue with navigation game'
$ws1.Range("C12").Value = 'negative'
$ws1.Range("D12").Value = 'In the navigation game, it was a bit disturbing for her and also for the patients I tried with today to collect the minerals when a note appears about maintaining a straight back - the correct posture window hides the route and is a bit stressful. Maybe after the first guidance, only voice guidance can be sufficient?'
$ws1.Range("A13").Value = ' '
$ws1.Range("A14").Value = 'feedback'
$ws1.Range("B14").Value = 'Amazing games and positive experience'
$ws1.Range("C14").Value = 'positive'
$ws1.Range("D14").Value = 'What a beauty!!!! Amazing games - both graphics, content and feedback that children receive! Yesterday I tried with my 5-year-old daughter, it was hard for her, she had to make an effort, but the satisfaction after the climb was worth the effort :).'
$ws1.Range("A15").Value = 'bug'
$ws1.Range("B15").Value = 'Issue with navigation game'
$ws1.Range("C15").Value = 'negative'
$ws1.Range("D15").Value = 'In the navigation game, it was a bit disturbing for her and also for the patients I tried with today to collect the minerals when a note appears about maintaining a straight back - the correct posture window hides the route and is a bit stressful. Maybe after the first guidance, only voice guidance can be sufficient?'

# Sheet3 data ('shitttt')
$ws3.Range("A1").Value = 'feedback'
$ws3.Range("B1").Value = 'Amazing games - graphics, content, and feedback'
$ws3.Range("C1").Value = 'positive'
$ws3.Range("D1").Value = 'What a beauty!!!! Amazing games - both graphics, content, and feedback children receive! Yesterday I tried it with my 5-year-old daughter, it was hard for her, she had to make an effort, but the satisfaction after climbing was worth the effort :).'
$ws3.Range("A2").Value = 'bug'
$ws3.Range("B2").Value = 'Navigation game issue with posture reminder'
$ws3.Range("C2").Value = 'negative'
$ws3.Range("D2").Value = 'In the navigation game, it disturbed her and also the patients I tried with today to collect the minerals when a notice appears about keeping a straight back - the correct posture window hides the path and is a bit stressful. Maybe after the first guidance, only a vocal guidance could be sufficient?'
$ws3.Range("A3").Value = 'feedback'
$ws3.Range("B3").Value = 'Amazing games, graphics, and feedback'
$ws3.Range("C3").Value = 'positive'
$ws3.Range("D3").Value = 'What a beauty!!!! Amazing games - both graphics, content, and feedback that children receive! Yesterday I tried it with my 5-year-old daughter, it was hard for her, she had to make an effort, but the satisfaction after the climb was worth the effort :).'
$ws3.Range("A4").Value = 'bug'
$ws3.Range("B4").Value = 'Navigation game issue'
$ws3.Range("C4").Value = 'negative'
$ws3.Range("D4").Value = 'In the navigation game, it was a bit disturbing for her and also for the patients I tried with today to collect minerals when a note about keeping a straight back appears - the correct posture window hides the path and is a bit stressful. Maybe after the first guidance, only voice guidance can be sufficient?'

# Sheet4 data ('shitttt2')
$ws4.Range("A1").Value = 'feedback'
$ws4.Range("B1").Value = 'Amazing games and graphics'
$ws4.Range("C1").Value = 'positive'
$ws4.Range("D1").Value = 'What a beauty! Amazing games - both graphics, content, and feedback kids get! Yesterday I tried it with my 5-year-old daughter, it was hard for her, she had to make an effort, but the satisfaction after the climb was worth the effort :).'
$ws4.Range("A2").Value = 'bug'
$ws4.Range("B2").Value = 'Navigation issue'
$ws4.Range("C2").Value = 'negative'
$ws4.Range("D2").Value = 'In the navigation game, it was a bit disturbing for her and also for the patients I tried with today to collect the minerals when a note about keeping a straight back appears - the proper posture window hides the track and is a bit stressful. Maybe after the first guidance, only vocal guidance could be sufficient?'
$ws4.Range("A3").Value = ' '
$ws4.Range("A4").Value = 'feedback'
$ws4.Range("B4").Value = 'Amazing games, graphics, content, and feedback'
$ws4.Range("C4").Value = 'positive'
$ws4.Range("D4").Value = 'What a beauty!!!! Amazing games - both graphics, content and also feedback that children receive! Yesterday I tried with my 5-year-old daughter, it was hard for her, she had to make an effort, but the satisfaction after the climb was worth the effort :).'
$ws4.Range("A5").Value = 'bug'
$ws4.Range("B5").Value = 'Navigation game issue with posture reminder'
$ws4.Range("C5").Value = 'negative'
$ws4.Range("D5").Value = 'In the navigation game, it was a bit disturbing for her and also for the patients I tried with today to collect the minerals when a note about maintaining a straight back appears - the correct posture window hides the path and is a bit stressful. Maybe after the first direction, only a vocal direction could be sufficient?'
$ws4.Range("A6").Value = ' '

# Restore the active selection on Sheet1 to A2 and make Sheet1 the active sheet again
$ws1.Activate()
[void]$ws1.Range("A2").Select()
